$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.630.75"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.514.23"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'586.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").Value = "'132.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "3.514.97"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'7.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "4.113.07"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "'27.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D17").Value = "3.511.99"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "64.604.58"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "'10.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "'14.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "'392.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "3.654.87"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").Value = "'74.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'0.0000110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'8.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").Value = "3.517.32"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "'171.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.816"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'26.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "'42.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "2.476.06"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "'0.913"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.82%  "
